$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("E15").Value = 369
$ws1.Range("L15").Value = 2255.48
$ws1.Range("M15").Value = 786.91
$ws1.Range("K25").Value = 1011.87
$ws1.Range("M32").Value = 6715.25
$ws1.Range("M33").Value = 1934.19
$ws1.Range("M44").Value = 4248.33

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F15").Value = 5570.7
$ws2.Range("F25").Value = 3334.1
$ws2.Range("F32").Value = 10362.68
$ws2.Range("F33").Value = 1934.19
$ws2.Range("F44").Value = 4248.33
$ws2.Range("F62").Value = 67557.35000000001

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D4").Value = 1745.31
$ws3.Range("E4").Value = -955.9299999999999
$ws3.Range("F4").Value = 2.21098837061998

$ws3.Range("D10").Value = 4225.87
$ws3.Range("E10").Value = 5690.13
$ws3.Range("F10").Value = 0.4261668011294877

$ws3.Range("D11").Value = 8735.85
$ws3.Range("E11").Value = 7412.15
$ws3.Range("F11").Value = 0.5409864998761457

$ws3.Range("D12").Value = 31512.49
$ws3.Range("E12").Value = 18794.51
$ws3.Range("F12").Value = 0.6264036813962272

$ws3.Range("D14").Value = 70523.29000000001
$ws3.Range("E14").Value = 27338.59766749098
$ws3.Range("F14").Value = 0.7206410144020483
